$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-157)
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03).
$ws.Range("C2:C157").Value = 45233
